# "adding menu to main prj"
# Append three new rows (102-104) to the sheet with variant spellings of
# "Tm Bax" that already exists earlier in the list (row 15 in the shared
# strings table), then move the active selection past the new data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Values typed in this exact order so the new shared-string entries land
# with "TM Bax" before "Tm Bax" in the shared string table (index 98/99),
# matching how the rows reference them (A102 -> "Tm Bax", A103 -> "TM Bax").
$ws.Range("A103").Value = "TM Bax"
$ws.Range("A102").Value = "Tm Bax"
$ws.Range("A104").Value = "TM BAX"

# Reflect the new selection/active cell just past the appended rows.
$ws.Range("A105").Select()
